# "Add multi target simulation data"
# The "edges" sheet gets a new row (i7 -> u, label "=", vsa FALSE) describing
# the equality import used when u_clip is FALSE, and the existing i8 -> u
# equality-import row's note is clarified to say it applies when u_clip is TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("edges")

# Make "edges" the active/selected sheet (matches workbookView/activeTab and
# the sheetView tabSelected flip in the diff).
$ws.Select()

# --- Update existing row 20 note text -------------------------------------
# "equality to import" -> "equality to import (if u_clip = TRUE)"
$ws.Range("E20").Value = "equality to import (if u_clip = TRUE)"

# --- Append new row 21: i7 -> u, "=", FALSE, "...(if u_clip = FALSE)" -----
$ws.Range("A21").Value = "i7"
$ws.Range("B21").Value = "u"

# Column C holds a literal "=" sign. A plain assignment would be parsed as
# the start of a formula by the Value setter, so build it as a formula that
# evaluates to the literal text, then flatten it down to a plain value with
# copy / paste-special so the stored cell is a literal (shared) string, not
# a formula - matching how the rest of the sheet stores this column.
$ws.Range("C21").Formula = "=""="""
$ws.Range("C21").Copy()
$ws.Range("C21").PasteSpecial(-4163)

# Column D holds the literal text "FALSE" (not the boolean value FALSE) -
# same trick, otherwise Excel auto-coerces the string into a true boolean.
$ws.Range("D21").Formula = "=""FALSE"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)

$ws.Range("E21").Value = "equality to import (if u_clip = FALSE)"

# Clear clipboard marching ants / leftover copy source.
$excel.CutCopyMode = 0

# --- Update the view/selection state --------------------------------------
$ws.Range("E22").Select()
